# Add RAG embeddings infrastructure - resume data template cleanup
#
# Logical changes to the Education sheet: the hard-coded start_date /
# end_date sample values are cleared out, and the Education sheet becomes
# the active/selected sheet (with C2 as the selected cell), replacing
# Profile as the previously-active sheet.

$wb = $excel.ActiveWorkbook

$education = $wb.Worksheets.Item("Education")

# Clear the sample start_date (C2) and end_date (D2) values on the
# Education sheet, leaving the cells blank.
$education.Range("C2:D2").ClearContents()

# Make Education the active sheet and select cell C2, which becomes the
# new active tab / selection for the workbook.
$education.Activate()
$education.Range("C2").Select()

Write-Host "Cleared Education C2:D2 and activated Education!C2"
